$wb = $excel.ActiveWorkbook

# Rename the "size" sheet to "thickthin"
$sheet = $wb.Worksheets.Item("size")
$sheet.Name = "thickthin"

# Make the renamed sheet the active/selected tab
$sheet.Activate()
